# "add toast message and fix bug"
# The toast-message part of the commit lives in the FE JS code (not part of
# this workbook). The parts of the commit that touch
# "FE/assets/Danh sach nhan vien.xlsx" are:
#   - the sheet got renamed to "Danh sách nhân viên"
#   - the two trailing empty columns (S:T) were deleted
#   - the remaining columns were narrowed (AutoFit-style column resize)
#   - row heights were set explicitly (title rows + wrapped data rows)
#   - the view moved (scrolled down a row, selection now on H9)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the sheet --------------------------------------------------
# A plain `.Name = ...` keeps sheetId/r:id untouched. The real edit's
# sheetId incremented (1 -> 2) while r:id stayed rId1, which is what you get
# from "duplicate the sheet, drop the original, keep the duplicate" - so
# reproduce that sequence (it also preserves every cell/style/hyperlink).
$ws.Copy($null, $ws)
$original = $wb.Worksheets.Item("Sheet1")
$original.Delete()
$renamed = $wb.Worksheets.Item("Sheet1 (2)")
$renamed.Name = "Danh sách nhân viên"
$renamed.Activate()
$ws = $renamed

# --- Drop the two trailing (empty) columns ------------------------------
$ws.Columns("S:T").Delete()

# --- Narrow the remaining columns (AutoFit-style) -----------------------
$ws.Columns(1).ColumnWidth = 7.6666666666667
$ws.Columns(2).ColumnWidth = 7.6666666666667
$ws.Columns(3).ColumnWidth = 12.6666666666667
$ws.Columns(4).ColumnWidth = 19.6666666666667
$ws.Columns(5).ColumnWidth = 11
$ws.Columns(7).ColumnWidth = 18.6666666666667
$ws.Columns(8).ColumnWidth = 14.5
$ws.Columns(9).ColumnWidth = 13.3333333333333
$ws.Columns(10).ColumnWidth = 10.6666666666667
$ws.Columns(11).ColumnWidth = 7.6666666666667
$ws.Columns(12).ColumnWidth = 7.6666666666667
$ws.Columns(13).ColumnWidth = 7.6666666666667
$ws.Columns(14).ColumnWidth = 10.1666666666667
$ws.Columns(15).ColumnWidth = 10.8333333333333
$ws.Columns(16).ColumnWidth = 16.6666666666667
$ws.Columns(17).ColumnWidth = 7.6666666666667
$ws.Columns(18).ColumnWidth = 15.3333333333333

# --- Row heights ----------------------------------------------------------
$ws.Rows(1).RowHeight = 25.5
$ws.Rows(2).RowHeight = 25.5
$ws.Rows(4).RowHeight = 45
$ws.Rows(5).RowHeight = 30
$ws.Rows(6).RowHeight = 45
$ws.Rows(7).RowHeight = 45
$ws.Rows(8).RowHeight = 45
$ws.Rows(9).RowHeight = 30
$ws.Rows(10).RowHeight = 30
$ws.Rows(11).RowHeight = 45
$ws.Rows(12).RowHeight = 45
$ws.Rows(13).RowHeight = 45

# --- View: scrolled down one row, H9 selected ---------------------------
$ws.Range("H9").Select()
$excel.ActiveWindow.ScrollRow = 2
